$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "TC_011"
$ws.Range("C14").Value = "standard_user"
$ws.Range("D14").Value = "secret_sauce"
$ws.Range("K14").Value = "lohi"

# Row 15
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "TC_012"
$ws.Range("C15").Value = "standard_user"
$ws.Range("D15").Value = "secret_sauce"

# New column header (first used after rows 14/15 body data, per original authoring order)
$ws.Range("L1").Value = "ItemName"

$ws.Range("L15").Value = "Sauce Labs Bike Light"

# Update selection to reflect final state
$ws.Range("L16").Select()
